# math_table.xlsx edit:
#  - Remove the "^2" / " 제곱" lookup row
#  - Remove the "^3" / " 세제곱" lookup row
#  - Remove the "sqrt" / " 루트" lookup row
#  - Change the pronunciation text for "rightarrow" from " 즉" to " 은 즉"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row with "^2" / " 제곱" is row 53, immediately followed by the
# "^3" / " 세제곱" row 54 - delete both in one go so the rest of the
# table (rows 55-100) shifts up by two.
$ws.Rows("53:54").Delete()

# After the shift above, the "sqrt" / " 루트" pair that used to live at
# row 69 is now at row 67 - remove it as well.
$ws.Rows("67").Delete()

# The "rightarrow" row (now row 68) keeps its symbol but its Korean
# pronunciation changes from " 즉" to " 은 즉".
$ws.Range("B68").Value = " 은 즉"

# Restore the cursor/selection position recorded in the saved workbook.
$ws.Range("D70").Select() | Out-Null
